$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data in columns D (Fecha) and M:T (Volumen..Kg/unidad) for rows
# 6, 9, 10, 11, 12 is cyclically rotated between those rows. Column A-C
# and E-L stay as-is. We set the new values explicitly per the target
# (after) state described by the diff.

# Row 6 (was row 10's data)
$ws.Range("D6").Value = [DateTime]"2021-04-21"
$ws.Range("M6").Value = 150
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17000
$ws.Range("Q6").Value = "`$/caja 15 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1133
$ws.Range("T6").Value = 15

# Row 9 (was row 11's data)
$ws.Range("D9").Value = [DateTime]"2021-05-18"
$ws.Range("M9").Value = 16
$ws.Range("N9").Value = 240000
$ws.Range("O9").Value = 250000
$ws.Range("P9").Value = 245000
$ws.Range("Q9").Value = "`$/bins (450 kilos)"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 544
$ws.Range("T9").Value = 450

# Row 10 (was row 6's data)
$ws.Range("D10").Value = [DateTime]"2021-03-25"
$ws.Range("M10").Value = 15
$ws.Range("N10").Value = 360000
$ws.Range("O10").Value = 360000
$ws.Range("P10").Value = 360000
$ws.Range("Q10").Value = "`$/bins (450 kilos)"
$ws.Range("R10").Value = "Provincia del Elquí"
$ws.Range("S10").Value = 800
$ws.Range("T10").Value = 450

# Row 11 (was row 12's data)
$ws.Range("D11").Value = [DateTime]"2021-03-30"
$ws.Range("M11").Value = 8
$ws.Range("N11").Value = 280000
$ws.Range("O11").Value = 300000
$ws.Range("P11").Value = 290000
$ws.Range("Q11").Value = "`$/bins (400 kilos)"
$ws.Range("R11").Value = "Provincia del Elquí"
$ws.Range("S11").Value = 725
$ws.Range("T11").Value = 400

# Row 12 (was row 9's data)
$ws.Range("D12").Value = [DateTime]"2021-05-04"
$ws.Range("M12").Value = 12
$ws.Range("N12").Value = 250000
$ws.Range("O12").Value = 260000
$ws.Range("P12").Value = 255000
$ws.Range("Q12").Value = "`$/bins (400 kilos)"
$ws.Range("R12").Value = "Provincia de Limarí"
$ws.Range("S12").Value = 638
$ws.Range("T12").Value = 400
